# Generate Report for Handoff
# Update status text and timestamps on the three sheets, and shrink the
# width of the "Status"/"Latest Handoff Datetime" style columns that were
# formatted for the older, longer status string.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn (E2) / de-de (F2) status + generate date (G2) ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-18 19:03:30"

# --- zh-cn sheet: Status (C2) + Latest Handoff Datetime (H2) ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-18 19:03:25"

# --- de-de sheet: Status (C2) + Latest Handoff Datetime (H2) ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-18 19:03:30"

# --- Column width adjustments (characters) ---
# Target stored width is 17.2159881591797 characters; the ColumnWidth
# setter snaps to whole-pixel (1/6-character) increments, so feed a value
# from the middle of the input bucket that rounds to the closest
# representable width (17.1666... == 103/6).
$overview.Range("E1").ColumnWidth = 16.3333333333333
$overview.Range("F1").ColumnWidth = 16.3333333333333
$zhcn.Range("C1").ColumnWidth = 16.3333333333333
$dede.Range("C1").ColumnWidth = 16.3333333333333
